$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "All other diseases (Residual)" row (row 15) entirely,
# shifting all subsequent rows up by one.
$ws.Rows.Item(15).Delete()
